$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("__data")

# Row 6 - Leviathan Lance
$ws.Range("D6").Value = "深渊利维坦矛"
$ws.Range("F6").Value = "weapon:20020002|relic:30050003"
$ws.Range("G6").Value = "光束伤害:+28|牵引力:+20"
$ws.Range("I6").Value = "icons/synergy/leviathan.png"

# Row 7 - Singularity Waltz
$ws.Range("D7").Value = "奇点圆舞"
$ws.Range("F7").Value = "weapon:20020001|relic:30050001"
$ws.Range("G7").Value = "弹速:+18|环轨:+1|暴击:+6"
$ws.Range("I7").Value = "icons/synergy/singularity.png"

# Row 8 - Seraph Tide
$ws.Range("D8").Value = "炽天潮汐"
$ws.Range("F8").Value = "relic:30050002|skill:70040002"
$ws.Range("G8").Value = "减速:+12%|护盾:+30|持续:+2"
$ws.Range("I8").Value = "icons/synergy/seraph_tide.png"

# Row 9 - Undertow Battery
$ws.Range("D9").Value = "逆潮蓄能阵"
$ws.Range("F9").Value = "weapon:20020003|skill:70030002"
$ws.Range("G9").Value = "爆发伤害:+24|碎片:+1"
$ws.Range("I9").Value = "icons/synergy/undertow_battery.png"
